$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PCB size changed for cutoff: every placement's "Mid Y" (column C) value
# shifts by the new board-edge offset. Write the updated coordinates
# (values taken from the refreshed pick-and-place export) directly.
$values = @{
    2  = 41.988
    3  = 41.988
    4  = 29.992000000000001
    5  = 37.473999999999997
    6  = 30.013000000000002
    7  = 37.542999999999999
    8  = 30.015000000000001
    9  = 37.508000000000003
    10 = 29.99
    11 = 37.542999999999999
    12 = 65.492000000000004
    13 = 29.576000000000001
    14 = 26.271999999999998
    15 = 57.902999999999999
    16 = 51.036999999999999
    17 = 58.180999999999997
    18 = 61.594000000000001
    19 = 64.274000000000001
    20 = 45.162999999999997
    21 = 41.036000000000001
    22 = 38.893000000000001
    23 = 62.625999999999998
    24 = 60.959000000000003
    25 = 57.704000000000001
    26 = 38.539000000000001
    27 = 42.146999999999998
    28 = 45.401000000000003
    29 = 57.250999999999998
    30 = 27.224
    31 = 57.902999999999999
    32 = 57.902999999999999
    33 = 31.481000000000002
    34 = 26.43
    35 = 26.43
    36 = 26.43
    37 = 26.43
    38 = 26.43
    39 = 26.43
    40 = 26.43
    41 = 26.431000000000001
    42 = 26.431000000000001
    43 = 26.431000000000001
    44 = 26.431000000000001
    45 = 26.431000000000001
    46 = 26.43
    47 = 26.43
    48 = 52.069000000000003
    49 = 49.448999999999998
    50 = 52.624000000000002
    51 = 48.814
    52 = 50.719000000000001
    53 = 33.847000000000001
    54 = 33.732999999999997
    55 = 34.029000000000003
    56 = 33.892000000000003
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}

# Update the view state: scroll the window down so row 13 is the top
# visible row, then select A2:E56 (mirrors the saved sheetView/selection).
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A2:E56").Select()
